$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 3
$ws.Range("G2").Value2 = 27.13486633333333
$ws.Range("H2").Value2 = 81.40459899999999
$ws.Range("I2").Value2 = 0.04747038381101173
$ws.Range("J2").Value2 = 0.05206311700485852
$ws.Range("K2").Value2 = 3
$ws.Range("M2").Value2 = 138.2190853333333
$ws.Range("N2").Value2 = 414.657256
$ws.Range("O2").Value2 = 0.2249223651785973
$ws.Range("P2").Value2 = 0.2476599003709697
$ws.Range("Q2").Value2 = 3750.556405235594
$ws.Range("R2").Value2 = 33755.00764712034
$ws.Range("S2").Value2 = 0.01067715100270855
$ws.Range("T2").Value2 = 0.0128939463704254

$ws.Range("E3").Value2 = 3
$ws.Range("G3").Value2 = 27.13486633333333
$ws.Range("H3").Value2 = 81.40459899999999
$ws.Range("I3").Value2 = 0.04747038381101173
$ws.Range("J3").Value2 = 0.05206311700485852
$ws.Range("K3").Value2 = 3
$ws.Range("M3").Value2 = 147.91433
$ws.Range("N3").Value2 = 443.74299
$ws.Range("O3").Value2 = 0.2406993279341593
$ws.Range("P3").Value2 = 0.2650317656414439
$ws.Range("Q3").Value2 = 4013.635573334556
$ws.Range("R3").Value2 = 36122.72016001101
$ws.Range("S3").Value2 = 0.01142608948008712
$ws.Range("T3").Value2 = 0.01379837982459473

$ws.Range("E4").Value2 = 3
$ws.Range("G4").Value2 = 27.13486633333333
$ws.Range("H4").Value2 = 81.40459899999999
$ws.Range("I4").Value2 = 0.04747038381101173
$ws.Range("J4").Value2 = 0.05206311700485852
$ws.Range("K4").Value2 = 3
$ws.Range("M4").Value2 = 74.27261733333333
$ws.Range("N4").Value2 = 222.817852
$ws.Range("O4").Value2 = 0.1208629960061633
$ws.Range("P4").Value2 = 0.1330811078998542
$ws.Range("Q4").Value2 = 2015.377543566816
$ws.Range("R4").Value2 = 18138.39789210134
$ws.Range("S4").Value2 = 0.005737412808961352
$ws.Range("T4").Value2 = 0.00692861729172631

$ws.Range("E5").Value2 = 3
$ws.Range("G5").Value2 = 27.13486633333333
$ws.Range("H5").Value2 = 81.40459899999999
$ws.Range("I5").Value2 = 0.04747038381101173
$ws.Range("J5").Value2 = 0.05206311700485852
$ws.Range("K5").Value2 = 3
$ws.Range("M5").Value2 = 84.85695366666668
$ws.Range("N5").Value2 = 254.570861
$ws.Range("O5").Value2 = 0.138086767645209
$ws.Range("P5").Value2 = 0.1520460408212704
$ws.Range("Q5").Value2 = 2302.58209519886
$ws.Range("R5").Value2 = 20723.23885678974
$ws.Range("S5").Value2 = 0.006555031859340069
$ws.Range("T5").Value2 = 0.007915990813403297

$ws.Range("E6").Value2 = 3
$ws.Range("G6").Value2 = 27.13486633333333
$ws.Range("H6").Value2 = 81.40459899999999
$ws.Range("I6").Value2 = 0.04747038381101173
$ws.Range("J6").Value2 = 0.05206311700485852
$ws.Range("K6").Value2 = 2
$ws.Range("M6").Value2 = 169.2560955
$ws.Range("N6").Value2 = 338.512191
$ws.Range("O6").Value2 = 0.275428543235871
$ws.Range("P6").Value2 = 0.2021811852664618
$ws.Range("Q6").Value2 = 4592.741527494401
$ws.Range("R6").Value2 = 27556.44916496641
$ws.Range("S6").Value2 = 0.01307469865991464
$ws.Range("T6").Value2 = 0.01052618270470878

$ws.Range("E7").Value2 = 3
$ws.Range("G7").Value2 = 207.121208
$ws.Range("H7").Value2 = 621.3636240000001
$ws.Range("I7").Value2 = 0.3623427924198875
$ws.Range("J7").Value2 = 0.3973992557702412
$ws.Range("K7").Value2 = 3
$ws.Range("M7").Value2 = 138.2190853333333
$ws.Range("N7").Value2 = 414.657256
$ws.Range("O7").Value2 = 0.2249223651785973
$ws.Range("P7").Value2 = 0.2476599003709697
$ws.Range("Q7").Value2 = 28628.10392289509
$ws.Range("R7").Value2 = 257652.9353060558
$ws.Range("S7").Value2 = 0.08149899787649864
$ws.Range("T7").Value2 = 0.09841986009155546

$ws.Range("E8").Value2 = 3
$ws.Range("G8").Value2 = 207.121208
$ws.Range("H8").Value2 = 621.3636240000001
$ws.Range("I8").Value2 = 0.3623427924198875
$ws.Range("J8").Value2 = 0.3973992557702412
$ws.Range("K8").Value2 = 3
$ws.Range("M8").Value2 = 147.91433
$ws.Range("N8").Value2 = 443.74299
$ws.Range("O8").Value2 = 0.2406993279341593
$ws.Range("P8").Value2 = 0.2650317656414439
$ws.Range("Q8").Value2 = 30636.19471011064
$ws.Range("R8").Value2 = 275725.7523909958
$ws.Range("S8").Value2 = 0.08721566661725352
$ws.Range("T8").Value2 = 0.1053234264213828

$ws.Range("E9").Value2 = 3
$ws.Range("G9").Value2 = 207.121208
$ws.Range("H9").Value2 = 621.3636240000001
$ws.Range("I9").Value2 = 0.3623427924198875
$ws.Range("J9").Value2 = 0.3973992557702412
$ws.Range("K9").Value2 = 3
$ws.Range("M9").Value2 = 74.27261733333333
$ws.Range("N9").Value2 = 222.817852
$ws.Range("O9").Value2 = 0.1208629960061633
$ws.Range("P9").Value2 = 0.1330811078998542
$ws.Range("Q9").Value2 = 15383.43422340174
$ws.Range("R9").Value2 = 138450.9080106156
$ws.Range("S9").Value2 = 0.04379383547310695
$ws.Range("T9").Value2 = 0.05288633323648122

$ws.Range("E10").Value2 = 3
$ws.Range("G10").Value2 = 207.121208
$ws.Range("H10").Value2 = 621.3636240000001
$ws.Range("I10").Value2 = 0.3623427924198875
$ws.Range("J10").Value2 = 0.3973992557702412
$ws.Range("K10").Value2 = 3
$ws.Range("M10").Value2 = 84.85695366666668
$ws.Range("N10").Value2 = 254.570861
$ws.Range("O10").Value2 = 0.138086767645209
$ws.Range("P10").Value2 = 0.1520460408212704
$ws.Range("Q10").Value2 = 17575.67475064004
$ws.Range("R10").Value2 = 158181.0727557603
$ws.Range("S10").Value2 = 0.05003474498480122
$ws.Range("T10").Value2 = 0.06042298346518458

$ws.Range("E11").Value2 = 3
$ws.Range("G11").Value2 = 207.121208
$ws.Range("H11").Value2 = 621.3636240000001
$ws.Range("I11").Value2 = 0.3623427924198875
$ws.Range("J11").Value2 = 0.3973992557702412
$ws.Range("K11").Value2 = 2
$ws.Range("M11").Value2 = 169.2560955
$ws.Range("N11").Value2 = 338.512191
$ws.Range("O11").Value2 = 0.275428543235871
$ws.Range("P11").Value2 = 0.2021811852664618
$ws.Range("Q11").Value2 = 35056.52696132337
$ws.Range("R11").Value2 = 210339.1617679402
$ws.Range("S11").Value2 = 0.09979954746822724
$ws.Range("T11").Value2 = 0.08034665255563717

$ws.Range("E12").Value2 = 3
$ws.Range("G12").Value2 = 84.750407
$ws.Range("H12").Value2 = 254.251221
$ws.Range("I12").Value2 = 0.1482643879283573
$ws.Range("J12").Value2 = 0.1626088848807073
$ws.Range("K12").Value2 = 3
$ws.Range("M12").Value2 = 138.2190853333333
$ws.Range("N12").Value2 = 414.657256
$ws.Range("O12").Value2 = 0.2249223651785973
$ws.Range("P12").Value2 = 0.2476599003709697
$ws.Range("Q12").Value2 = 11714.12373716773
$ws.Range("R12").Value2 = 105427.1136345096
$ws.Range("S12").Value2 = 0.0333479768046032
$ws.Range("T12").Value2 = 0.04027170022899046

$ws.Range("E13").Value2 = 3
$ws.Range("G13").Value2 = 84.750407
$ws.Range("H13").Value2 = 254.251221
$ws.Range("I13").Value2 = 0.1482643879283573
$ws.Range("J13").Value2 = 0.1626088848807073
$ws.Range("K13").Value2 = 3
$ws.Range("M13").Value2 = 147.91433
$ws.Range("N13").Value2 = 443.74299
$ws.Range("O13").Value2 = 0.2406993279341593
$ws.Range("P13").Value2 = 0.2650317656414439
$ws.Range("Q13").Value2 = 12535.79966863231
$ws.Range("R13").Value2 = 112822.1970176908
$ws.Range("S13").Value2 = 0.03568713853092508
$ws.Range("T13").Value2 = 0.04309651986892015

$ws.Range("E14").Value2 = 3
$ws.Range("G14").Value2 = 84.750407
$ws.Range("H14").Value2 = 254.251221
$ws.Range("I14").Value2 = 0.1482643879283573
$ws.Range("J14").Value2 = 0.1626088848807073
$ws.Range("K14").Value2 = 3
$ws.Range("M14").Value2 = 74.27261733333333
$ws.Range("N14").Value2 = 222.817852
$ws.Range("O14").Value2 = 0.1208629960061633
$ws.Range("P14").Value2 = 0.1330811078998542
$ws.Range("Q14").Value2 = 6294.634547955254
$ws.Range("R14").Value2 = 56651.71093159728
$ws.Range("S14").Value2 = 0.0179196781260413
$ws.Range("T14").Value2 = 0.02164017055428438

$ws.Range("E15").Value2 = 3
$ws.Range("G15").Value2 = 84.750407
$ws.Range("H15").Value2 = 254.251221
$ws.Range("I15").Value2 = 0.1482643879283573
$ws.Range("J15").Value2 = 0.1626088848807073
$ws.Range("K15").Value2 = 3
$ws.Range("M15").Value2 = 84.85695366666668
$ws.Range("N15").Value2 = 254.570861
$ws.Range("O15").Value2 = 0.138086767645209
$ws.Range("P15").Value2 = 0.1520460408212704
$ws.Range("Q15").Value2 = 7191.661360030143
$ws.Range("R15").Value2 = 64724.95224027129
$ws.Range("S15").Value2 = 0.0204733500859222
$ws.Range("T15").Value2 = 0.02472403714847329

$ws.Range("E16").Value2 = 3
$ws.Range("G16").Value2 = 84.750407
$ws.Range("H16").Value2 = 254.251221
$ws.Range("I16").Value2 = 0.1482643879283573
$ws.Range("J16").Value2 = 0.1626088848807073
$ws.Range("K16").Value2 = 2
$ws.Range("M16").Value2 = 169.2560955
$ws.Range("N16").Value2 = 338.512191
$ws.Range("O16").Value2 = 0.275428543235871
$ws.Range("P16").Value2 = 0.2021811852664618
$ws.Range("Q16").Value2 = 14344.52298085587
$ws.Range("R16").Value2 = 86067.13788513522
$ws.Range("S16").Value2 = 0.04083624438086551
$ws.Range("T16").Value2 = 0.03287645708003904

$ws.Range("E17").Value2 = 3
$ws.Range("G17").Value2 = 101.3352343333333
$ws.Range("H17").Value2 = 304.005703
$ws.Range("I17").Value2 = 0.1772782813185584
$ws.Range("J17").Value2 = 0.1944298562963656
$ws.Range("K17").Value2 = 3
$ws.Range("M17").Value2 = 138.2190853333333
$ws.Range("N17").Value2 = 414.657256
$ws.Range("O17").Value2 = 0.2249223651785973
$ws.Range("P17").Value2 = 0.2476599003709697
$ws.Range("Q17").Value2 = 14006.46340159233
$ws.Range("R17").Value2 = 126058.170614331
$ws.Range("S17").Value2 = 0.0398738503289669
$ws.Range("T17").Value2 = 0.04815247883949986

$ws.Range("E18").Value2 = 3
$ws.Range("G18").Value2 = 101.3352343333333
$ws.Range("H18").Value2 = 304.005703
$ws.Range("I18").Value2 = 0.1772782813185584
$ws.Range("J18").Value2 = 0.1944298562963656
$ws.Range("K18").Value2 = 3
$ws.Range("M18").Value2 = 147.91433
$ws.Range("N18").Value2 = 443.74299
$ws.Range("O18").Value2 = 0.2406993279341593
$ws.Range("P18").Value2 = 0.2650317656414439
$ws.Range("Q18").Value2 = 14988.93329180799
$ws.Range("R18").Value2 = 134900.3996262719
$ws.Range("S18").Value2 = 0.04267076317069984
$ws.Range("T18").Value2 = 0.05153008810763798

$ws.Range("E19").Value2 = 3
$ws.Range("G19").Value2 = 101.3352343333333
$ws.Range("H19").Value2 = 304.005703
$ws.Range("I19").Value2 = 0.1772782813185584
$ws.Range("J19").Value2 = 0.1944298562963656
$ws.Range("K19").Value2 = 3
$ws.Range("M19").Value2 = 74.27261733333333
$ws.Range("N19").Value2 = 222.817852
$ws.Range("O19").Value2 = 0.1208629960061633
$ws.Range("P19").Value2 = 0.1330811078998542
$ws.Range("Q19").Value2 = 7526.433082023328
$ws.Range("R19").Value2 = 67737.89773820995
$ws.Range("S19").Value2 = 0.02142638420698443
$ws.Range("T19").Value2 = 0.02587494068472978

$ws.Range("E20").Value2 = 3
$ws.Range("G20").Value2 = 101.3352343333333
$ws.Range("H20").Value2 = 304.005703
$ws.Range("I20").Value2 = 0.1772782813185584
$ws.Range("J20").Value2 = 0.1944298562963656
$ws.Range("K20").Value2 = 3
$ws.Range("M20").Value2 = 84.85695366666668
$ws.Range("N20").Value2 = 254.570861
$ws.Range("O20").Value2 = 0.138086767645209
$ws.Range("P20").Value2 = 0.1520460408212704
$ws.Range("Q20").Value2 = 8598.999284624477
$ws.Range("R20").Value2 = 77390.99356162029
$ws.Range("S20").Value2 = 0.02447978484097777
$ws.Range("T20").Value2 = 0.02956228986731095

$ws.Range("E21").Value2 = 3
$ws.Range("G21").Value2 = 101.3352343333333
$ws.Range("H21").Value2 = 304.005703
$ws.Range("I21").Value2 = 0.1772782813185584
$ws.Range("J21").Value2 = 0.1944298562963656
$ws.Range("K21").Value2 = 2
$ws.Range("M21").Value2 = 169.2560955
$ws.Range("N21").Value2 = 338.512191
$ws.Range("O21").Value2 = 0.275428543235871
$ws.Range("P21").Value2 = 0.2021811852664618
$ws.Range("Q21").Value2 = 17151.60609983755
$ws.Range("R21").Value2 = 102909.6365990253
$ws.Range("S21").Value2 = 0.04882749877092947
$ws.Range("T21").Value2 = 0.03931005879718703

$ws.Range("E22").Value2 = 2
$ws.Range("G22").Value2 = 151.2750305
$ws.Range("H22").Value2 = 302.550061
$ws.Range("I22").Value2 = 0.2646441545221851
$ws.Range("J22").Value2 = 0.1934988860478274
$ws.Range("K22").Value2 = 3
$ws.Range("M22").Value2 = 138.2190853333333
$ws.Range("N22").Value2 = 414.657256
$ws.Range("O22").Value2 = 0.2249223651785973
$ws.Range("P22").Value2 = 0.2476599003709697
$ws.Range("Q22").Value2 = 20909.0963494821
$ws.Range("R22").Value2 = 125454.5780968926
$ws.Range("S22").Value2 = 0.05952438916582006
$ws.Range("T22").Value2 = 0.04792191484049855

$ws.Range("E23").Value2 = 2
$ws.Range("G23").Value2 = 151.2750305
$ws.Range("H23").Value2 = 302.550061
$ws.Range("I23").Value2 = 0.2646441545221851
$ws.Range("J23").Value2 = 0.1934988860478274
$ws.Range("K23").Value2 = 3
$ws.Range("M23").Value2 = 147.91433
$ws.Range("N23").Value2 = 443.74299
$ws.Range("O23").Value2 = 0.2406993279341593
$ws.Range("P23").Value2 = 0.2650317656414439
$ws.Range("Q23").Value2 = 22375.74478213706
$ws.Range("R23").Value2 = 134254.4686928224
$ws.Range("S23").Value2 = 0.06369967013519376
$ws.Range("T23").Value2 = 0.05128335141890823

$ws.Range("E24").Value2 = 2
$ws.Range("G24").Value2 = 151.2750305
$ws.Range("H24").Value2 = 302.550061
$ws.Range("I24").Value2 = 0.2646441545221851
$ws.Range("J24").Value2 = 0.1934988860478274
$ws.Range("K24").Value2 = 3
$ws.Range("M24").Value2 = 74.27261733333333
$ws.Range("N24").Value2 = 222.817852
$ws.Range("O24").Value2 = 0.1208629960061633
$ws.Range("P24").Value2 = 0.1330811078998542
$ws.Range("Q24").Value2 = 11235.59245241483
$ws.Range("R24").Value2 = 67413.55471448896
$ws.Range("S24").Value2 = 0.03198568539106933
$ws.Range("T24").Value2 = 0.02575104613263251

$ws.Range("E25").Value2 = 2
$ws.Range("G25").Value2 = 151.2750305
$ws.Range("H25").Value2 = 302.550061
$ws.Range("I25").Value2 = 0.2646441545221851
$ws.Range("J25").Value2 = 0.1934988860478274
$ws.Range("K25").Value2 = 3
$ws.Range("M25").Value2 = 84.85695366666668
$ws.Range("N25").Value2 = 254.570861
$ws.Range("O25").Value2 = 0.138086767645209
$ws.Range("P25").Value2 = 0.1520460408212704
$ws.Range("Q25").Value2 = 12836.73825406209
$ws.Range("R25").Value2 = 77020.42952437252
$ws.Range("S25").Value2 = 0.03654385587416777
$ws.Range("T25").Value2 = 0.02942073952689831

$ws.Range("E26").Value2 = 2
$ws.Range("G26").Value2 = 151.2750305
$ws.Range("H26").Value2 = 302.550061
$ws.Range("I26").Value2 = 0.2646441545221851
$ws.Range("J26").Value2 = 0.1934988860478274
$ws.Range("K26").Value2 = 2
$ws.Range("M26").Value2 = 169.2560955
$ws.Range("N26").Value2 = 338.512191
$ws.Range("O26").Value2 = 0.275428543235871
$ws.Range("P26").Value2 = 0.2021811852664618
$ws.Range("Q26").Value2 = 25604.22100907341
$ws.Range("R26").Value2 = 102416.8840362936
$ws.Range("S26").Value2 = 0.07289055395593419
$ws.Range("T26").Value2 = 0.03912183412888976
